$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(34).Insert()
$ws.Cells.Item(34,1).Value = 182
$ws.Cells.Item(34,2).Value = "Tr. bright orange"
$ws.Cells.Item(34,4).Value = 225
$ws.Cells.Item(34,5).Value = 141
$ws.Cells.Item(34,6).Value = 10
$ws.Cells.Item(34,7).Value = "Generic"
$ws.Cells.Item(34,8).Value = $true
$ws.Cells.Item(34,9).Value = $true
$ws.Cells.Item(34,10).Value = $false
$ws.Cells.Item(34,11).Value = $false
